$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column F (reuse the header formatting from E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Add time_taken values for each data row
$ws.Range("F2").Value = "2021-10-05 10:51:38.190492"
$ws.Range("F3").Value = "2021-10-05 10:51:38.190501"
$ws.Range("F4").Value = "2021-10-05 10:51:38.190505"
$ws.Range("F5").Value = "2021-10-05 10:51:38.190507"
$ws.Range("F6").Value = "2021-10-05 10:51:38.190510"
$ws.Range("F7").Value = "2021-10-05 10:51:38.190513"
$ws.Range("F8").Value = "2021-10-05 10:51:38.190515"
$ws.Range("F9").Value = "2021-10-05 10:51:38.190518"
$ws.Range("F10").Value = "2021-10-05 10:51:38.190520"
$ws.Range("F11").Value = "2021-10-05 10:51:38.190523"
$ws.Range("F12").Value = "2021-10-05 10:51:38.190525"
$ws.Range("F13").Value = "2021-10-05 10:51:38.190528"
$ws.Range("F14").Value = "2021-10-05 10:51:38.190530"
$ws.Range("F15").Value = "2021-10-05 10:51:38.190533"
$ws.Range("F16").Value = "2021-10-05 10:51:38.190535"
$ws.Range("F17").Value = "2021-10-05 10:51:38.190538"
